$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select the entire first row (matches the target's post-insert selection
# state) and insert a new blank row above it, pushing all existing rows
# down by one.
$ws.Rows.Item(1).Select()
$ws.Rows.Item(1).Insert()
